$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BOM formulas to account for new back panel screws and nuts
$ws.Range("B19").Formula = "=62+6"
$ws.Range("B23").Formula = "=58+6"

# Update the view state: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B24").Select()
